$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.633.16"
$ws.Range("D3").Value = "1.592.60"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "210.45"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.816.06"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.586.87"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "26.612.51"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "206.57"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "145.59"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "7.17"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "0.660"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "1.278.30"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "1.49"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "63.33"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "0.918"
$ws.Range("E45").Value = "  +9.30%  "
$ws.Range("D46").Value = "1.728.18"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "89.59"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.04%  "
